$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add "Poverty rate(%)" data (column K) for years 2000-2023 (rows 5-28)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$povertyRates = @{
    5  = 8.8
    6  = 10.4
    7  = 10.1
    8  = 10.9
    9  = 10.6
    10 = 9.8
    11 = 9.3
    12 = 8
    13 = 9.1
    14 = 10.4
    15 = 10.7
    16 = 12
    17 = 11.6
    18 = 10.8
    19 = 11.4
    20 = 10.6
    21 = 9.3
    22 = 9.5
    23 = 8.8
    24 = 9.3
    25 = 9.3
    26 = 11.2
    27 = 10.2
    28 = 10.1
}

foreach ($row in $povertyRates.Keys) {
    $ws1.Cells.Item($row, 11).Value = $povertyRates[$row]
}

# ---------------------------------------------------------------------------
# Sheet2: add a new Indicator / source row documenting the new "Poverty rate"
# column, citing the Statista source used for the data above.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A10").Value = "Poverty rate"
$ws2.Range("B10").Value = "American Community Survey (ACS) obtained from Statista website: https://www.statista.com/statistics/205456/poverty-rate-in-hawaii/"

# Match the look of the rows above (centered, same font) by copying the
# formatting from the previous row's label cell.
$ws2.Range("A9").Copy()
$ws2.Range("A10").PasteSpecial(-4122)

# Restore selections: Sheet2's own active cell moves to A11, then Sheet1
# (the originally active sheet) becomes active again with K14 selected.
$ws2.Range("A11").Select()
$ws1.Range("K14").Select()
